$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(32, 1).Value = '(홍보)[한국항공우주연구원] 초격차 스타트업 인재발굴 및 채용지원 프로그램 참여자 모집 홍보'
$ws.Cells.Item(32, 2).Value = '관리자'
$ws.Cells.Item(32, 3).Value = '''2025-11-05'
$ws.Cells.Item(32, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=137'

$ws.Cells.Item(33, 1).Value = '(홍보)IITP 2026 AI·ICT 산업·기술전망 컨퍼런스'
$ws.Cells.Item(33, 2).Value = '관리자'
$ws.Cells.Item(33, 3).Value = '''2025-10-29'
$ws.Cells.Item(33, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=136'

$ws.Cells.Item(34, 1).Value = '(안내)2025년 디지털혁신네트워크 지역인재 채용 면접 희망자 신청 안내(~25.10.31.16:00)'
$ws.Cells.Item(34, 2).Value = '관리자'
$ws.Cells.Item(34, 3).Value = '''2025-10-27'
$ws.Cells.Item(34, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=135'

$ws.Cells.Item(35, 1).Value = '한국언론진흥재단 빅카인즈 OPEN API 사용 신청 안내'
$ws.Cells.Item(35, 2).Value = '관리자'
$ws.Cells.Item(35, 3).Value = '''2025-10-24'
$ws.Cells.Item(35, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=134'

$ws.Cells.Item(36, 1).Value = '국립공주대학교 SW중심대학사업단-한국언론진흥재단 빅카인즈 OPEN API 활용을 위한 업무협약 체결'
$ws.Cells.Item(36, 2).Value = '관리자'
$ws.Cells.Item(36, 3).Value = '''2025-10-24'
$ws.Cells.Item(36, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=133'

$ws.Cells.Item(37, 1).Value = '(공지) 2025년 SW알고리즘 경진대회 참여신청 안내 (구글폼 참여신청:11.7. (금) 까지)'
$ws.Cells.Item(37, 2).Value = '관리자'
$ws.Cells.Item(37, 3).Value = '''2025-10-24'
$ws.Cells.Item(37, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=132'

$ws.Cells.Item(38, 1).Value = '2025년 SW전문가 특강(마음Ai, 11월5일(수), 구글폼사전신청 10월 29일(수)까지) 안내'
$ws.Cells.Item(38, 2).Value = '관리자'
$ws.Cells.Item(38, 3).Value = '''2025-10-22'
$ws.Cells.Item(38, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=131'

$ws.Cells.Item(39, 1).Value = '(홍보)[경기대학교] 2025 SW전공교수 TOPCIT 릴레이 온라인 특강 안내'
$ws.Cells.Item(39, 2).Value = '관리자'
$ws.Cells.Item(39, 3).Value = '''2025-10-22'
$ws.Cells.Item(39, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=129'

$ws.Cells.Item(40, 1).Value = '(채용) 국립공주대학교 SW중심대학사업단 산학협력초빙교수(비전임) 채용 재공고'
$ws.Cells.Item(40, 2).Value = '관리자'
$ws.Cells.Item(40, 3).Value = '''2025-10-20'
$ws.Cells.Item(40, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=128'

$ws.Cells.Item(41, 1).Value = '(홍보) 온라인 특강 안내( 25.10.21.(화) 19시~20시, 온라인, 삼성전자 박수홍 오픈소스그룹장)'
$ws.Cells.Item(41, 2).Value = '관리자'
$ws.Cells.Item(41, 3).Value = '''2025-10-14'
$ws.Cells.Item(41, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=125'

$ws.Cells.Item(42, 1).Value = '(공지-보완) 프로그래머스 & N사 현직자와 함께하는 백엔드 취업 집중 과정 설명회'
$ws.Cells.Item(42, 2).Value = '관리자'
$ws.Cells.Item(42, 3).Value = '''2025-10-01'
$ws.Cells.Item(42, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=124'

$ws.Cells.Item(43, 1).Value = '(홍보) 2025 데이터안심구역 활용 공동 경진대회( ~ 25.10.10. 18시)'
$ws.Cells.Item(43, 2).Value = '관리자'
$ws.Cells.Item(43, 3).Value = '''2025-09-24'
$ws.Cells.Item(43, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=123'

$ws.Cells.Item(44, 1).Value = '2025년 SW전문가 특강(4차) 안내'
$ws.Cells.Item(44, 2).Value = '관리자'
$ws.Cells.Item(44, 3).Value = '''2025-09-24'
$ws.Cells.Item(44, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=122'

$ws.Cells.Item(45, 1).Value = '2025 SW중심대학 에세이 공모전 (''25. 10. 1.(수) ~ 10.14.(화), 17시) 접수 안내'
$ws.Cells.Item(45, 2).Value = '관리자'
$ws.Cells.Item(45, 3).Value = '''2025-09-22'
$ws.Cells.Item(45, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=121'

$ws.Cells.Item(46, 1).Value = '(마감) 2025 SW인재페스티벌 우수작품 경진대회 출품작 모집 안내'
$ws.Cells.Item(46, 2).Value = '관리자'
$ws.Cells.Item(46, 3).Value = '''2025-09-22'
$ws.Cells.Item(46, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=120'

$ws.Cells.Item(47, 1).Value = '(공지) 2025학년도 2학기 산학캡스톤디자인 프로젝트 운영 안내'
$ws.Cells.Item(47, 2).Value = '관리자'
$ws.Cells.Item(47, 3).Value = '''2025-09-22'
$ws.Cells.Item(47, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=119'

$ws.Cells.Item(48, 1).Value = '(신청마감) [제24회 TOPCIT 정기평가 시행 안내] -구글폼신청[9/15(월) 08:30~9/16(화)15시까지 (선착순)]'
$ws.Cells.Item(48, 2).Value = '관리자'
$ws.Cells.Item(48, 3).Value = '''2025-09-08'
$ws.Cells.Item(48, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=118'

$ws.Cells.Item(49, 1).Value = '2025년 COEIC 집중교육 교육대상자 선정 안내(SW마일리지, value-up 마일리지 부여)'
$ws.Cells.Item(49, 2).Value = '관리자'
$ws.Cells.Item(49, 3).Value = '''2025-08-18'
$ws.Cells.Item(49, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=109'

$ws.Cells.Item(50, 1).Value = '(마감) 2025년 충남 학생 정보과학 챌린지 개최'
$ws.Cells.Item(50, 2).Value = '관리자'
$ws.Cells.Item(50, 3).Value = '''2025-08-12'
$ws.Cells.Item(50, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=106'

$ws.Cells.Item(51, 1).Value = '(마감) 2025년 COEIC 집중교육 참여 신청 안내(SW마일리지 부여/~25. 8. 14. 13:00)'
$ws.Cells.Item(51, 2).Value = '관리자'
$ws.Cells.Item(51, 3).Value = '''2025-08-12'
$ws.Cells.Item(51, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=105'

$ws.Cells.Item(52, 1).Value = '(마감) 2025년 TOPCIT 파일럿테스트 참여 신청 안내(SW마일리지 50점 & 기프티콘 ) - 선착순25명+대기5명'
$ws.Cells.Item(52, 2).Value = '관리자'
$ws.Cells.Item(52, 3).Value = '''2025-08-07'
$ws.Cells.Item(52, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=104'

$ws.Cells.Item(53, 1).Value = '(마감)2025학년도 2학기 「ICT 학점연계 프로젝트 인턴십 사업」 학생 모집 안내'
$ws.Cells.Item(53, 2).Value = '관리자'
$ws.Cells.Item(53, 3).Value = '''2025-07-30'
$ws.Cells.Item(53, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=103'

$ws.Cells.Item(54, 1).Value = '지역 SW전공 학생 취업희망 수요조사 실시_충남테크노파크, 지역SW중심대학사업단'
$ws.Cells.Item(54, 2).Value = '관리자'
$ws.Cells.Item(54, 3).Value = '''2025-07-30'
$ws.Cells.Item(54, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=102'

$ws.Cells.Item(55, 1).Value = '(안내)2025년 충남 인공지능 실증랩 운영계획 안내_충남연구원'
$ws.Cells.Item(55, 2).Value = '관리자'
$ws.Cells.Item(55, 3).Value = '''2025-07-28'
$ws.Cells.Item(55, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=101'

$ws.Cells.Item(56, 1).Value = '2025년 2학기 창업교과목 수강신청 안내(SW마일리지 50점)'
$ws.Cells.Item(56, 2).Value = '관리자'
$ws.Cells.Item(56, 3).Value = '''2025-07-25'
$ws.Cells.Item(56, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=100'

$ws.Cells.Item(57, 1).Value = '2025년도 2학기 주관학과 SW마일리지 점수표 공지(2025.07.10. 기준)'
$ws.Cells.Item(57, 2).Value = '관리자'
$ws.Cells.Item(57, 3).Value = '''2025-07-18'
$ws.Cells.Item(57, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=98'

$ws.Cells.Item(58, 1).Value = '(홍보) 충남콘텐츠진흥원 「피지컬 AI 대학생 창업 경진대회」 참여 안내(~7.20.까지)'
$ws.Cells.Item(58, 2).Value = '관리자'
$ws.Cells.Item(58, 3).Value = '''2025-07-18'
$ws.Cells.Item(58, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=97'

$ws.Cells.Item(59, 1).Value = '2025년 주관학과 대상 SW교육(COEIC) 진단평가 및 교육 신청 안내(~25.7.20.)'
$ws.Cells.Item(59, 2).Value = '관리자'
$ws.Cells.Item(59, 3).Value = '''2025-07-14'
$ws.Cells.Item(59, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=95'

$ws.Cells.Item(60, 1).Value = '제23회 TOPCIT 정기평가 성적우수자 총장상 시상식 개최'
$ws.Cells.Item(60, 2).Value = '관리자'
$ws.Cells.Item(60, 3).Value = '''2025-07-10'
$ws.Cells.Item(60, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=94'

$ws.Cells.Item(61, 1).Value = '25년 1학기 SW중심대학사업 SW마일리지 장학금 지급 명단 안내 (총 140명)'
$ws.Cells.Item(61, 2).Value = '관리자'
$ws.Cells.Item(61, 3).Value = '''2025-07-10'
$ws.Cells.Item(61, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=93'

$ws.Cells.Item(62, 1).Value = '제23회 TOPCIT 정기평가 성적우수자 시상 안내'
$ws.Cells.Item(62, 2).Value = '관리자'
$ws.Cells.Item(62, 3).Value = '''2025-07-08'
$ws.Cells.Item(62, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=92'

$ws.Cells.Item(63, 1).Value = '2025년 1학기 SW마일리지 적립내역 확인 및 추가 신청 안내( ~25.7.8.(화) 17:00 )'
$ws.Cells.Item(63, 2).Value = '관리자'
$ws.Cells.Item(63, 3).Value = '''2025-07-07'
$ws.Cells.Item(63, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=91'

$ws.Cells.Item(64, 1).Value = '2025년 SW전문가 특강(3차) 안내'
$ws.Cells.Item(64, 2).Value = '관리자'
$ws.Cells.Item(64, 3).Value = '''2025-06-17'
$ws.Cells.Item(64, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=84'

$ws.Cells.Item(65, 1).Value = '2025 SW중심대학 디지털 경진대회 참가자 모집 및 신청 안내'
$ws.Cells.Item(65, 2).Value = '관리자'
$ws.Cells.Item(65, 3).Value = '''2025-06-04'
$ws.Cells.Item(65, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=82'

$ws.Cells.Item(66, 1).Value = '(홍보) LG Aimers 7기 모집 (접수기간 : ~25. 6. 19.(목))'
$ws.Cells.Item(66, 2).Value = '관리자'
$ws.Cells.Item(66, 3).Value = '''2025-05-27'
$ws.Cells.Item(66, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=81'

$ws.Cells.Item(67, 1).Value = '2025 글로벌 SW교육 프로그램(SPP) 해외교육 프로그램 최종합격자 발표'
$ws.Cells.Item(67, 2).Value = '관리자'
$ws.Cells.Item(67, 3).Value = '''2025-05-13'
$ws.Cells.Item(67, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=80'

$ws.Cells.Item(68, 1).Value = '[경기대학교] 2025 AI-Powered SW전공교수 TOPCIT 릴레이 온라인 특강 안내'
$ws.Cells.Item(68, 2).Value = '관리자'
$ws.Cells.Item(68, 3).Value = '''2025-05-12'
$ws.Cells.Item(68, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=79'

$ws.Cells.Item(69, 1).Value = '2025 주관학과 대상 SW집중교육 운영 안내(PCCP 자격취득 대비용)'
$ws.Cells.Item(69, 2).Value = '관리자'
$ws.Cells.Item(69, 3).Value = '''2025-05-09'
$ws.Cells.Item(69, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=78'

$ws.Cells.Item(70, 1).Value = '2025 글로벌 SW교육 프로그램(SPP) 참가 신청 안내(25.7.6 - 8.2, 4주간, 미국 LA(USC) 등/해외교육비 지원)'
$ws.Cells.Item(70, 2).Value = '관리자'
$ws.Cells.Item(70, 3).Value = '''2025-05-08'
$ws.Cells.Item(70, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=69'

$ws.Cells.Item(71, 1).Value = '2025년 SW전문가 특강 (2차) 안내'
$ws.Cells.Item(71, 2).Value = '관리자'
$ws.Cells.Item(71, 3).Value = '''2025-05-07'
$ws.Cells.Item(71, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=67'

$ws.Cells.Item(72, 1).Value = '2025 PCCE 합격자 명단 및 PCCP 오프라인 교육 대상자 안내'
$ws.Cells.Item(72, 2).Value = '관리자'
$ws.Cells.Item(72, 3).Value = '''2025-05-02'
$ws.Cells.Item(72, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=65'

$ws.Cells.Item(73, 1).Value = '2025 PCCE(코딩역량인증시험) 안내'
$ws.Cells.Item(73, 2).Value = '관리자'
$ws.Cells.Item(73, 3).Value = '''2025-04-23'
$ws.Cells.Item(73, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=64'

$ws.Cells.Item(74, 1).Value = '2025년 WCRC 물류로봇경진대회 개최'
$ws.Cells.Item(74, 2).Value = '관리자'
$ws.Cells.Item(74, 3).Value = '''2025-04-15'
$ws.Cells.Item(74, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=60'

$ws.Cells.Item(75, 1).Value = '[ 제23회 TOPCIT 정기평가 시행 안내 ] - 구글폼 신청 ~4.16일까지'
$ws.Cells.Item(75, 2).Value = '관리자'
$ws.Cells.Item(75, 3).Value = '''2025-04-11'
$ws.Cells.Item(75, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=59'

$ws.Cells.Item(76, 1).Value = '2025 ICCAS 해외교육 프로그램 최종합격자 발표'
$ws.Cells.Item(76, 2).Value = '관리자'
$ws.Cells.Item(76, 3).Value = '''2025-04-10'
$ws.Cells.Item(76, 4).Value = 'https://swknu.kongju.ac.kr/community/noticedetail.do?seq=58'
